$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 370.8
$ws.Range("I2").Value = 370.8
$ws.Range("K2").Value = 370.8
$ws.Range("M2").Value = -257.8
$ws.Range("H33").Value = 201.6842
$ws.Range("I33").Value = 165
$ws.Range("K33").Value = 165
$ws.Range("M33").Value = 64
$ws.Range("H70").Value = 5699.25
$ws.Range("I70").Value = 6898.5
$ws.Range("K70").Value = 20695.5
$ws.Range("M70").Value = -20425.5
$ws.Range("H73").Value = 5699.25
$ws.Range("I73").Value = 6898.5
$ws.Range("K73").Value = 20695.5
$ws.Range("M73").Value = -19759.5
$ws.Range("H138").Value = 256634.27
$ws.Range("J138").Value = 316357.88
$ws.Range("L138").Value = 949073.64
$ws.Range("N138").Value = -959353.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 10003
$ws.Range("I13").Value = 10003
$ws.Range("K13").Value = 10003
$ws.Range("M13").Value = -9859
$ws.Range("H32").Value = 2148.08
$ws.Range("I32").Value = 2129.3738
$ws.Range("K32").Value = 2129.3738
$ws.Range("M32").Value = -1842.3738
$ws.Range("H45").Value = 28999.117
$ws.Range("I45").Value = 62447.145
$ws.Range("J45").Value = 5585.5
$ws.Range("K45").Value = 62447.145
$ws.Range("L45").Value = 5585.5
$ws.Range("M45").Value = -62070.145
$ws.Range("N45").Value = -6339.5
$ws.Range("H61").Value = 6372.5264
$ws.Range("I61").Value = 5643.4546
$ws.Range("K61").Value = 5643.4546
$ws.Range("M61").Value = -5431.4546
$ws.Range("H74").Value = 509344.72
$ws.Range("I74").Value = 621476.9
$ws.Range("K74").Value = 621476.9
$ws.Range("M74").Value = -620602.9
$ws.Range("H77").Value = 509344.72
$ws.Range("I77").Value = 621476.9
$ws.Range("K77").Value = 3107384.5
$ws.Range("M77").Value = -3103016.5
$ws.Range("H80").Value = 84955
$ws.Range("J80").Value = 84955
$ws.Range("L80").Value = 84955
$ws.Range("N80").Value = -86951
$ws.Range("H83").Value = 84955
$ws.Range("J83").Value = 84955
$ws.Range("L83").Value = 254865
$ws.Range("N83").Value = -264849
$ws.Range("H110").Value = 2321.0386
$ws.Range("I110").Value = 1353.2858
$ws.Range("J110").Value = 6385.6
$ws.Range("K110").Value = 1353.2858
$ws.Range("L110").Value = 6385.6
$ws.Range("M110").Value = 691.7141999999999
$ws.Range("N110").Value = -10475.6
$ws.Range("H136").Value = 6372.5264
$ws.Range("I136").Value = 5643.4546
$ws.Range("K136").Value = 16930.3638
$ws.Range("M136").Value = -14380.3638
$ws.Range("H139").Value = 86198.8
$ws.Range("J139").Value = 86198.8
$ws.Range("L139").Value = 86198.8
$ws.Range("N139").Value = -96478.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6571.375
$ws.Range("I86").Value = 6674.2
$ws.Range("J86").Value = 6400
$ws.Range("K86").Value = 6674.2
$ws.Range("L86").Value = 6400
$ws.Range("M86").Value = -5551.2
$ws.Range("N86").Value = -8646
$ws.Range("H89").Value = 6571.375
$ws.Range("I89").Value = 6674.2
$ws.Range("J89").Value = 6400
$ws.Range("K89").Value = 33371
$ws.Range("L89").Value = 32000
$ws.Range("M89").Value = -27755
$ws.Range("N89").Value = -43232
$ws.Range("H107").Value = 1433.579
$ws.Range("I107").Value = 1412.7333
$ws.Range("J107").Value = 1511.75
$ws.Range("K107").Value = 1412.7333
$ws.Range("L107").Value = 1511.75
$ws.Range("M107").Value = 507.2666999999999
$ws.Range("N107").Value = -5351.75
$ws.Range("H134").Value = 2505.84
$ws.Range("I134").Value = 1821.6875
$ws.Range("K134").Value = 5465.0625
$ws.Range("M134").Value = -2930.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4032.2556
$ws.Range("J31").Value = 5469.7334
$ws.Range("L31").Value = 5469.7334
$ws.Range("N31").Value = -6059.7334
$ws.Range("H34").Value = 4032.2556
$ws.Range("J34").Value = 5469.7334
$ws.Range("L34").Value = 5469.7334
$ws.Range("N34").Value = -5873.7334
$ws.Range("H96").Value = 10000
$ws.Range("J96").Value = 10000
$ws.Range("L96").Value = 10000
$ws.Range("N96").Value = -15492
$ws.Range("H99").Value = 5864.8125
$ws.Range("J99").Value = 5699.7144
$ws.Range("L99").Value = 5699.7144
$ws.Range("N99").Value = -8695.714400000001
$ws.Range("H106").Value = 252329.67
$ws.Range("J106").Value = 252329.67
$ws.Range("L106").Value = 252329.67
$ws.Range("N106").Value = -254853.67
$ws.Range("H126").Value = 5864.8125
$ws.Range("J126").Value = 5699.7144
$ws.Range("L126").Value = 17099.1432
$ws.Range("N126").Value = -22039.1432
$ws.Range("H134").Value = 2921.8
$ws.Range("I134").Value = 2007.3793
$ws.Range("J134").Value = 7341.5
$ws.Range("K134").Value = 6022.1379
$ws.Range("L134").Value = 22024.5
$ws.Range("M134").Value = -3487.1379
$ws.Range("N134").Value = -27094.5
$ws.Range("H141").Value = 268331.5
$ws.Range("J141").Value = 268331.5
$ws.Range("L141").Value = 268331.5
$ws.Range("N141").Value = -278691.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 62822220
$ws.Range("J4").Value = 6235999.5
$ws.Range("L4").Value = 18707998.5
$ws.Range("N4").Value = -18708222.5
$ws.Range("H5").Value = 1944.8148
$ws.Range("I5").Value = 417.72223
$ws.Range("K5").Value = 1253.16669
$ws.Range("M5").Value = -1141.16669
$ws.Range("H80").Value = 7898.8
$ws.Range("J80").Value = 7749
$ws.Range("L80").Value = 23247
$ws.Range("N80").Value = -25119
$ws.Range("H83").Value = 7898.8
$ws.Range("J83").Value = 7749
$ws.Range("L83").Value = 69741
$ws.Range("N83").Value = -79101
$ws.Range("H113").Value = 2054.25
$ws.Range("J113").Value = 2431.889
$ws.Range("L113").Value = 7295.667
$ws.Range("N113").Value = -11635.667
$ws.Range("H122").Value = 1452.909
$ws.Range("J122").Value = 1444.6666
$ws.Range("L122").Value = 13001.9994
$ws.Range("N122").Value = -17901.9994
$ws.Range("H135").Value = 1944.8148
$ws.Range("I135").Value = 417.72223
$ws.Range("K135").Value = 3759.50007
$ws.Range("M135").Value = -1224.50007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 2165
$ws.Range("I13").Value = 249.5
$ws.Range("J13").Value = 5996
$ws.Range("K13").Value = 249.5
$ws.Range("L13").Value = 5996
$ws.Range("M13").Value = -110.5
$ws.Range("N13").Value = -6274
$ws.Range("H126").Value = 14121.333
$ws.Range("I126").Value = 5907
$ws.Range("K126").Value = 17721
$ws.Range("M126").Value = -15251
$ws.Range("H139").Value = 99917.60000000001
$ws.Range("J139").Value = 99917.60000000001
$ws.Range("L139").Value = 99917.60000000001
$ws.Range("N139").Value = -110197.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6301.08
$ws.Range("I7").Value = 4650.7856
$ws.Range("J7").Value = 8401.454
$ws.Range("K7").Value = 4650.7856
$ws.Range("L7").Value = 8401.454
$ws.Range("M7").Value = -4538.7856
$ws.Range("N7").Value = -8625.454
$ws.Range("H30").Value = 875
$ws.Range("J30").Value = 1000
$ws.Range("L30").Value = 1000
$ws.Range("N30").Value = -1216
$ws.Range("H35").Value = 5391
$ws.Range("J35").Value = 1000
$ws.Range("L35").Value = 1000
$ws.Range("N35").Value = -1672
$ws.Range("H98").Value = 300000
$ws.Range("J98").Value = 300000
$ws.Range("L98").Value = 300000
$ws.Range("N98").Value = -305990
$ws.Range("H100").Value = 6540.3335
$ws.Range("I100").Value = 5048.5
$ws.Range("K100").Value = 5048.5
$ws.Range("M100").Value = -4507.5
$ws.Range("H126").Value = 6301.08
$ws.Range("I126").Value = 4650.7856
$ws.Range("J126").Value = 8401.454
$ws.Range("K126").Value = 13952.3568
$ws.Range("L126").Value = 25204.362
$ws.Range("M126").Value = -11482.3568
$ws.Range("N126").Value = -30144.362
$ws.Range("H132").Value = 3584.3103
$ws.Range("I132").Value = 2550.3157
$ws.Range("J132").Value = 5548.9
$ws.Range("K132").Value = 7650.9471
$ws.Range("L132").Value = 16646.7
$ws.Range("M132").Value = -5120.9471
$ws.Range("N132").Value = -21706.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 4750
$ws.Range("I12").Value = 4750
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 4750
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -4608
$ws.Range("N12").ClearContents()
$ws.Range("H113").Value = 1202.4
$ws.Range("I113").Value = 1633.5
$ws.Range("J113").Value = 915
$ws.Range("K113").Value = 4900.5
$ws.Range("L113").Value = 2745
$ws.Range("M113").Value = -2730.5
$ws.Range("N113").Value = -7085
$ws.Range("H122").Value = 12503256
$ws.Range("I122").Value = 2793.2856
$ws.Range("J122").Value = 41671000
$ws.Range("K122").Value = 8379.856800000001
$ws.Range("L122").Value = 125013000
$ws.Range("M122").Value = -5929.856800000001
$ws.Range("N122").Value = -125017900
$ws.Range("H126").Value = 2254.1667
$ws.Range("I126").Value = 2254.1667
$ws.Range("K126").Value = 6762.500100000001
$ws.Range("M126").Value = -4292.500100000001
